# Expenses Details.xlsx update
# 1] Added "Battery pack delivery" Dunzo Bill (row 16->15 numbered, new row 17)
# 2] Added "Charger Delivery ... + 4 bulbs for load experiment" expense (new row 18)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# Row 16 gets its "Sr. No" filled in (it was blank before)
$ws.Range("A16").Value = 15

# --- Row 17: Battery pack delivery (Dunzo bill) ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 44368
$ws.Range("C17").Value = "Battery pack delivery"
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 93
$ws.Range("F17").Value = "Akshay"

# --- Row 18: Charger delivery + bulbs for load experiment ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 44374
$ws.Range("C18").Value = "Charger Delivery from SSK home (Petrol 50Rs) + 4 bulbs for load experiment (15W 15 Rs each)"
$ws.Range("E18").Value = 110
$ws.Range("F18").Value = "Akshay"

# Bill link for the new Dunzo delivery bill (added last so the shared-string
# order matches: Battery pack delivery, Charger Delivery..., Bill_15)
$ws.Hyperlinks.Add($ws.Range("G17"), "Bills\Bill_15.pdf", "", "", "Bill_15") | Out-Null

# Reflect the final cursor/scroll position left behind by the edit session
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("H22").Select() | Out-Null
